$wb = $excel.ActiveWorkbook

# Sheet "展览" (Exhibition) - update "想去人数" (want-to-go count) values
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F4").Value = 143
$ws1.Range("F5").Value = 3007
$ws1.Range("F6").Value = 304
$ws1.Range("F7").Value = 405

# Sheet "全部类型" (All types) - update "想去人数" (want-to-go count) values
$ws2 = $wb.Worksheets.Item("全部类型")
$ws2.Range("F4").Value = 143
$ws2.Range("F5").Value = 3007
$ws2.Range("F6").Value = 304
$ws2.Range("F9").Value = 405
